# feat: add 2022-Q1 data
#
# The workbook's existing "总计" (grand-total) sheet is duplicated in place;
# the original becomes the new "2022-Q1" per-fund holdings sheet (reusing
# its sheetId/underlying part, as a renamed worksheet does), and the
# duplicate becomes the new "总计" sheet, carrying the old totals table plus
# a freshly prepended 2022-Q1 summary row. Using Worksheet.Copy (rather than
# Worksheets.Add) means both resulting sheets keep the original sheetPr /
# formatting metadata instead of starting from a blank sheet.

$wb = $excel.ActiveWorkbook

$total    = $wb.Worksheets.Item("总计")
$totalIdx = $total.Index

$total.Copy($null, $total)
$newTotal = $wb.Worksheets.Item($totalIdx + 1)

$total.Name    = "2022-Q1"
$newTotal.Name = "总计"

$q1 = $total

# ---------------------------------------------------------------------
# 1. "2022-Q1" sheet: per-fund holdings table (was the old 4-column
#    totals table; grows to 8 columns / 18 rows).
# ---------------------------------------------------------------------

# extend the bold/centered/bordered header style (currently only B1:D1)
# across the new E1:H1 header cells
$q1.Range("D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)

# extend the index-column style (currently only A2:A6) down through A18
$q1.Range("A6").Copy()
$q1.Range("A7:A18").PasteSpecial(-4122)

$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# code, name, fund size, total stock position, position share, held value (亿元), position rank
$fundRows = @(
    @("008545", "泓德丰润三年持有期混合",   "84.75", "91.10", "4.13", "3.5002", 10),
    @("010864", "泓德卓远混合A",            "39.08", "91.88", "4.82", "1.8837", 5),
    @("005395", "泓德臻远回报灵活配置混合", "33.94", "93.62", "5.24", "1.7785", 6),
    @("001500", "泓德远见回报混合",          "26.71", "93.72", "6.55", "1.7495", 7),
    @("004965", "泓德致远混合A",            "21.03", "46.32", "4.77", "1.0031", 6),
    @("010865", "泓德卓远混合C",            "12.25", "91.88", "4.82", "0.5904", 5),
    @("004966", "泓德致远混合C",            "2.92",  "46.32", "4.77", "0.1393", 6),
    @("001628", "招商体育文化休闲股票",      "2.95",  "83.21", "3.78", "0.1115", 8),
    @("159855", "银华中证影视主题ETF",       "0.96",  "97.27", "6.40", "0.0614", 5),
    @("290012", "泰信行业精选灵活配置混合A", "0.76",  "92.62", "7.44", "0.0565", 1),
    @("001223", "鹏华文化传媒娱乐股票",      "0.91",  "83.63", "4.54", "0.0413", 3),
    @("516620", "国泰中证影视主题ETF",       "0.33",  "96.08", "6.46", "0.0213", 5),
    @("006227", "华宝科技先锋混合A",         "0.46",  "91.73", "2.44", "0.0112", 9),
    @("008112", "中泰中证500指数增强A",      "0.61",  "92.46", "0.86", "0.0052", 9),
    @("008113", "中泰中证500指数增强C",      "0.46",  "92.46", "0.86", "0.0040", 9),
    @("010842", "华宝科技先锋混合C",         "0.04",  "91.73", "2.44", "0.0010", 9),
    @("002583", "泰信行业精选灵活配置混合C", "0.00",  "92.62", "7.44", "0",      1)
)

$r = 2
foreach ($row in $fundRows) {
    # numeric-looking D/E/F/G columns are stored as TEXT in the source
    # workbook (matching every other quarter sheet); a leading "'" keeps
    # Excel from re-interpreting them as numbers, except for the one true
    # zero market value (G18), which is a genuine numeric 0.
    $q1.Cells.Item($r, 1).Value = $r - 2
    $q1.Cells.Item($r, 2).Value = "'" + $row[0]
    $q1.Cells.Item($r, 3).Value = $row[1]
    $q1.Cells.Item($r, 4).Value = "'" + $row[2]
    $q1.Cells.Item($r, 5).Value = "'" + $row[3]
    $q1.Cells.Item($r, 6).Value = "'" + $row[4]
    if ($row[5] -eq "0") {
        $q1.Cells.Item($r, 7).Value = 0
    } else {
        $q1.Cells.Item($r, 7).Value = "'" + $row[5]
    }
    $q1.Cells.Item($r, 8).Value = $row[6]
    $r++
}

# ---------------------------------------------------------------------
# 2. "总计" sheet: quarterly summary table, now starting with 2022-Q1.
# ---------------------------------------------------------------------

# extend the index-column style (currently only A2:A6) down to the new A7
$newTotal.Range("A6").Copy()
$newTotal.Range("A7").PasteSpecial(-4122)

$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

$totalRows = @(
    @("2022-Q1", 17, 10.96),
    @("2021-Q4", 41, 23.31),
    @("2021-Q3", 19, 14.28),
    @("2021-Q2", 16, 12.01),
    @("2021-Q1", 23, 23.07),
    @("2020-Q4", 14, 7.62)
)

$r = 2
foreach ($row in $totalRows) {
    $newTotal.Cells.Item($r, 1).Value = $r - 2
    $newTotal.Cells.Item($r, 2).Value = $row[0]
    $newTotal.Cells.Item($r, 3).Value = $row[1]
    $newTotal.Cells.Item($r, 4).Value = $row[2]
    $r++
}
